$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 80
$ws.Range("H80").Value = 3558.9678
$ws.Range("I80").Value = 383.875
$ws.Range("J80").Value = 6945.7334
$ws.Range("K80").Value = 1151.625
$ws.Range("L80").Value = 20837.2002
$ws.Range("M80").Value = -153.625
$ws.Range("N80").Value = -22833.2002

# Row 83
$ws.Range("H83").Value = 3558.9678
$ws.Range("I83").Value = 383.875
$ws.Range("J83").Value = 6945.7334
$ws.Range("K83").Value = 3454.875
$ws.Range("L83").Value = 62511.6006
$ws.Range("M83").Value = 1537.125
$ws.Range("N83").Value = -72495.60060000001

# Row 131
$ws.Range("H131").Value = 5887.5
$ws.Range("I131").Value = 1850
$ws.Range("J131").Value = 7233.3335
$ws.Range("K131").Value = 5550
$ws.Range("L131").Value = 21700.0005
$ws.Range("M131").Value = -510
$ws.Range("N131").Value = -31780.0005

# Row 132
$ws.Range("H132").Value = 5130088
$ws.Range("I132").Value = 1252.9246
$ws.Range("J132").Value = 27782444
$ws.Range("K132").Value = 3758.7738
$ws.Range("L132").Value = 83347332
$ws.Range("M132").Value = -1228.7738
$ws.Range("N132").Value = -83352392

# Row 134
$ws.Range("H134").Value = 55683.25
$ws.Range("J134").Value = 55683.25
$ws.Range("L134").Value = 55683.25
$ws.Range("N134").Value = -65823.25

# Row 137
$ws.Range("H137").Value = 1623.6511
$ws.Range("I137").Value = 1266.2368
$ws.Range("J137").Value = 4340
$ws.Range("K137").Value = 3798.7104
$ws.Range("L137").Value = 13020
$ws.Range("M137").Value = -1248.7104
$ws.Range("N137").Value = -18120

# Row 138
$ws.Range("H138").Value = 2442.8208
$ws.Range("I138").Value = 1343.9487
$ws.Range("J138").Value = 3973.3928
$ws.Range("K138").Value = 4031.8461
$ws.Range("L138").Value = 11920.1784
$ws.Range("M138").Value = 1108.1539
$ws.Range("N138").Value = -22200.1784

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 118905.7
$ws.Range("I2").Value = 167865.92
$ws.Range("J2").Value = 1401.2
$ws.Range("K2").Value = 167865.92
$ws.Range("L2").Value = 1401.2
$ws.Range("M2").Value = -167752.92
$ws.Range("N2").Value = -1627.2

# Row 32
$ws.Range("H32").Value = 12023.466
$ws.Range("I32").Value = 8458.014999999999
$ws.Range("K32").Value = 8458.014999999999
$ws.Range("M32").Value = -8171.014999999999

# Row 116
$ws.Range("H116").Value = 118905.7
$ws.Range("I116").Value = 167865.92
$ws.Range("J116").Value = 1401.2
$ws.Range("K116").Value = 167865.92
$ws.Range("L116").Value = 1401.2
$ws.Range("M116").Value = -165571.92
$ws.Range("N116").Value = -5989.2

# Row 130
$ws.Range("H130").Value = 40000
$ws.Range("J130").Value = 40000
$ws.Range("L130").Value = 40000
$ws.Range("N130").Value = -50040

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 118905.7
$ws.Range("I3").Value = 167865.92
$ws.Range("J3").Value = 1401.2
$ws.Range("K3").Value = 167865.92
$ws.Range("L3").Value = 1401.2
$ws.Range("M3").Value = -167751.92
$ws.Range("N3").Value = -1629.2

# Row 81
$ws.Range("H81").Value = 28673.334
$ws.Range("J81").Value = 28673.334
$ws.Range("L81").Value = 28673.334
$ws.Range("N81").Value = -30795.334

# Row 84
$ws.Range("H84").Value = 28673.334
$ws.Range("J84").Value = 28673.334
$ws.Range("L84").Value = 86020.00199999999
$ws.Range("N84").Value = -96628.00199999999

# Row 122
$ws.Range("H122").Value = 67678.336
$ws.Range("J122").Value = 67678.336
$ws.Range("L122").Value = 67678.336
$ws.Range("N122").Value = -77478.336

# Row 133
$ws.Range("H133").Value = 42243.332
$ws.Range("J133").Value = 42692
$ws.Range("L133").Value = 42692
$ws.Range("N133").Value = -52812

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 5522659.5
$ws.Range("I58").Value = 9805580
$ws.Range("J58").Value = 668683.4399999999
$ws.Range("K58").Value = 9805580
$ws.Range("L58").Value = 668683.4399999999
$ws.Range("M58").Value = -9805377
$ws.Range("N58").Value = -669089.4399999999

# Row 106
$ws.Range("H106").Value = 40000
$ws.Range("J106").Value = 40000
$ws.Range("L106").Value = 40000
$ws.Range("N106").Value = -42524

# Row 136
$ws.Range("H136").Value = 5522659.5
$ws.Range("I136").Value = 9805580
$ws.Range("J136").Value = 668683.4399999999
$ws.Range("K136").Value = 29416740
$ws.Range("L136").Value = 2006050.32
$ws.Range("M136").Value = -29414190
$ws.Range("N136").Value = -2011150.32

$ws = $wb.Worksheets.Item("CUL")
# Row 101
$ws.Range("H101").Value = 9533.333000000001
$ws.Range("J101").Value = 9533.333000000001
$ws.Range("L101").Value = 28599.999
$ws.Range("N101").Value = -33467.999

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2850.6428
$ws.Range("I7").Value = 2100.4
$ws.Range("J7").Value = 4726.25
$ws.Range("K7").Value = 2100.4
$ws.Range("L7").Value = 4726.25
$ws.Range("M7").Value = -1988.4
$ws.Range("N7").Value = -4950.25

# Row 46
$ws.Range("H46").Value = 995
$ws.Range("I46").Value = 647
$ws.Range("J46").Value = 1169
$ws.Range("K46").Value = 647
$ws.Range("L46").Value = 1169
$ws.Range("M46").Value = -459
$ws.Range("N46").Value = -1545

# Row 61
$ws.Range("H61").Value = 2205.2942
$ws.Range("I61").Value = 2006.1666
$ws.Range("J61").Value = 2683.2
$ws.Range("K61").Value = 2006.1666
$ws.Range("L61").Value = 2683.2
$ws.Range("M61").Value = -1804.1666
$ws.Range("N61").Value = -3087.2

# Row 113
$ws.Range("H113").Value = 2205.2942
$ws.Range("I113").Value = 2006.1666
$ws.Range("J113").Value = 2683.2
$ws.Range("K113").Value = 2006.1666
$ws.Range("L113").Value = 2683.2
$ws.Range("M113").Value = 163.8334
$ws.Range("N113").Value = -7023.2

# Row 126
$ws.Range("H126").Value = 2850.6428
$ws.Range("I126").Value = 2100.4
$ws.Range("J126").Value = 4726.25
$ws.Range("K126").Value = 6301.200000000001
$ws.Range("L126").Value = 14178.75
$ws.Range("M126").Value = -3831.200000000001
$ws.Range("N126").Value = -19118.75

$ws = $wb.Worksheets.Item("WVR")
# Row 135
$ws.Range("H135").Value = 40081.168
$ws.Range("J135").Value = 40081.168
$ws.Range("L135").Value = 40081.168
$ws.Range("N135").Value = -50221.168

# Row 136
$ws.Range("H136").Value = 3404734.8
$ws.Range("I136").Value = 4046.2632
$ws.Range("J136").Value = 5558504
$ws.Range("K136").Value = 12138.7896
$ws.Range("L136").Value = 16675512
$ws.Range("M136").Value = -9588.7896
$ws.Range("N136").Value = -16680612
